# Saldo.xlsx update
#
# The underlying "Saldo" values for two accounts (TIAGO) were corrected
# upward, which moves their rows to a new position because the sheet is
# kept sorted by Saldo (column C) descending:
#   - account 004911541 (TIAGO): 184.13   -> 75184.13
#   - account 005924958 (TIAGO): 84.24    -> 25084.24
# Three other accounts were removed entirely:
#   - 001879977 (THAISSA) 10010.11
#   - 002064834 (RAFAELA)  3230.03
#   - 004357848 (AURELIO)  1796.62
#
# Implemented as: delete the five affected rows (the two old TIAGO rows plus
# the three removed accounts), then re-insert the two TIAGO rows with their
# corrected balances at their new, sorted position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Assert-Account($rowNum, $expectedAccount) {
    $actual = $ws.Cells.Item($rowNum, 1).Value2
    if ($actual -ne $expectedAccount) {
        throw "Row $rowNum expected account '$expectedAccount' but found '$actual'"
    }
}

# --- 1) Remove the old rows -------------------------------------------------
# Delete from the bottom up so earlier row numbers stay valid while we work.

Assert-Account 172 "005924958"
$ws.Rows.Item(172).Delete()   # 005924958 TIAGO   84.24  (old position)

Assert-Account 100 "004911541"
$ws.Rows.Item(100).Delete()   # 004911541 TIAGO  184.13  (old position)

Assert-Account 12 "004357848"
$ws.Rows.Item(12).Delete()    # 004357848 AURELIO 1796.62 (removed account)

Assert-Account 10 "002064834"
$ws.Rows.Item(10).Delete()    # 002064834 RAFAELA 3230.03 (removed account)

Assert-Account 7 "001879977"
$ws.Rows.Item(7).Delete()     # 001879977 THAISSA 10010.11 (removed account)

# --- 2) Re-insert the corrected TIAGO rows in sorted (descending) order ----

# 005924958 / TIAGO / 25084.24 belongs between 008035153 (30051.49) and
# 005142611 (22434.22), i.e. right before row 5 (GUILHERME).
Assert-Account 5 "005142611"
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "005924958"
$ws.Cells.Item(5, 2).Value = "TIAGO"
$ws.Cells.Item(5, 3).Value = 25084.24

# 004911541 / TIAGO / 75184.13 is now the largest balance after the header
# row, so it belongs right before row 2 (ANA / 008012870).
Assert-Account 2 "008012870"
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "004911541"
$ws.Cells.Item(2, 2).Value = "TIAGO"
$ws.Cells.Item(2, 3).Value = 75184.13
